# Generate Report for Handoff
#
# Rows 4-7 (the "Ready for handoff" files) in both the "zh-cn" and
# "de-de" localization-status sheets get refreshed by the handoff-report
# generator:
#   - Priority moves from the placeholder "low" to "ht" (matching the
#     already-handed-off rows 2-3).
#   - Latest Handoff Datetime is bumped to the new generation timestamp.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$zhcnHandoffTime = "2016-08-23 06:29:27"
$dedeHandoffTime = "2016-08-23 06:29:33"

foreach ($row in 4..7) {
    $zhcn.Range("E$row").Value = "ht"
    $zhcn.Range("H$row").Value = $zhcnHandoffTime

    $dede.Range("E$row").Value = "ht"
    $dede.Range("H$row").Value = $dedeHandoffTime
}
